$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1:J255")
$rng.AutoFilter(4, @("MS"), 7)
